# [Feat 3695] mostly functionnal, but not out of the wood yet.
#
# Insert a new "TC_STEP_CALL_DATASET" column into the STEPS sheet, right
# before the existing "TC_STEP_ACTION" column (old column G), shifting the
# remaining call-step columns one slot to the right. The new column is left
# blank for every data row; only the header is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STEPS")

# Insert a blank column at G - existing G:K (TC_STEP_ACTION .. TC_STEP_CUF_<CODE>)
# shift right to H:L, inheriting their row formatting automatically.
$ws.Columns("G").Insert()

# New header for the freshly inserted column.
$ws.Range("G1").Value = "TC_STEP_CALL_DATASET"

# Move the active selection to G2, like the author left it after typing the
# new header and dropping down into the first data row.
$ws.Range("G2").Select() | Out-Null
